$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel sending clusters with updated TPM-derived cluster names:
#   row 2: "ECs"  -> "FAPs"
#   row 3: "FAPs" -> "MuSCs"
$ws.Range("A2").Value = "FAPs"
$ws.Range("A3").Value = "MuSCs"

# Updated numeric results recomputed from the new TPM values (row 2)
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.344207
$ws.Range("H2").Value = 1.032621
$ws.Range("I2").Value = 0.6985282229833164
$ws.Range("J2").Value = 0.6985282229833165
$ws.Range("M2").Value = 0.1642713333333333
$ws.Range("Q2").Value = 0.05654334283266666
$ws.Range("R2").Value = 0.5088900854939999
$ws.Range("S2").Value = 0.6985282229833164
$ws.Range("T2").Value = 0.6985282229833165

# Updated numeric results recomputed from the new TPM values (row 3)
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.1485533333333333
$ws.Range("H3").Value = 0.44566
$ws.Range("I3").Value = 0.3014717770166836
$ws.Range("J3").Value = 0.3014717770166836
$ws.Range("M3").Value = 0.1642713333333333
$ws.Range("Q3").Value = 0.02440305413777778
$ws.Range("R3").Value = 0.21962748724
$ws.Range("S3").Value = 0.3014717770166836
$ws.Range("T3").Value = 0.3014717770166836
